# cv121121a.xlsx — "correção nos dados e inicio da analise PNAD 2009"
#
# Row 6 ("grandes regiões e unidades da federação") was a stray label row
# with no data underneath it. Removing it entirely shifts every following
# region row (norte, rondônia, acre, ... distrito federal) up by one row,
# so each region keeps its own data while the dimension shrinks from
# A1:H38 to A1:H37 and the now-unused shared string is dropped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
